$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws "D2" "61.342.60"
$ws.Range("E2").Value = "  -4.06%  "

Set-TextCell $ws "D3" "2.997.76"
$ws.Range("E3").Value = "  -2.96%  "

$ws.Range("E4").Value = "  +0.07%  "

Set-TextCell $ws "D5" "548.73"
$ws.Range("E5").Value = "  +1.00%  "

Set-TextCell $ws "D6" "133.73"
$ws.Range("E6").Value = "  -4.21%  "

$ws.Range("E7").Value = "  +0.02%  "

Set-TextCell $ws "D8" "2.993.34"
$ws.Range("E8").Value = "  -2.96%  "

Set-TextCell $ws "D9" "0.496"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  -5.33%  "

$ws.Range("E11").Value = "  -7.91%  "

Set-TextCell $ws "D12" "0.450"
$ws.Range("E12").Value = "  -1.70%  "

Set-TextCell $ws "D13" "34.41"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("E14").Value = "  -2.58%  "

Set-TextCell $ws "D15" "3.491.28"
$ws.Range("E15").Value = "  -2.85%  "

Set-TextCell $ws "D16" "61.539.36"
$ws.Range("E16").Value = "  -3.78%  "

$ws.Range("E17").Value = "  -2.30%  "

Set-TextCell $ws "D18" "3.001.81"
$ws.Range("E18").Value = "  -3.02%  "

Set-TextCell $ws "D19" "6.66"
$ws.Range("E19").Value = "  +0.02%  "

Set-TextCell $ws "D20" "471.96"
$ws.Range("E20").Value = "  -1.84%  "

Set-TextCell $ws "D21" "13.26"
$ws.Range("E21").Value = "  -1.09%  "

Set-TextCell $ws "D22" "0.673"
$ws.Range("E22").Value = "  -3.70%  "

Set-TextCell $ws "D23" "7.03"
$ws.Range("E23").Value = "  -1.17%  "

Set-TextCell $ws "D24" "80.03"
$ws.Range("E24").Value = "  +1.10%  "

Set-TextCell $ws "D25" "12.07"
$ws.Range("E25").Value = "  -2.49%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  -0.16%  "

Set-TextCell $ws "D28" "7.79"
$ws.Range("E28").Value = "  -3.26%  "

Set-TextCell $ws "D29" "1.00"
$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  +0.14%  "

Set-TextCell $ws "D31" "25.61"
$ws.Range("E31").Value = "  -2.62%  "

Set-TextCell $ws "D32" "1.13"
$ws.Range("E32").Value = "  -2.49%  "

Set-TextCell $ws "D33" "5.52"
$ws.Range("E33").Value = "  +2.55%  "

Set-TextCell $ws "D34" "2.29"
$ws.Range("E34").Value = "  -2.79%  "

Set-TextCell $ws "D35" "55.26"
$ws.Range("E35").Value = "  -3.67%  "

Set-TextCell $ws "D36" "5.89"
$ws.Range("E36").Value = "  -2.13%  "

Set-TextCell $ws "D37" "452.26"
$ws.Range("E37").Value = "  -8.50%  "

Set-TextCell $ws "D38" "3.184.95"
$ws.Range("E38").Value = "  -2.60%  "

Set-TextCell $ws "D39" "0.0798"
$ws.Range("E39").Value = "  -0.08%  "

$ws.Range("E40").Value = "  -5.32%  "

Set-TextCell $ws "D41" "0.117"
$ws.Range("E41").Value = "  -2.39%  "

Set-TextCell $ws "D42" "8.14"
$ws.Range("E42").Value = "  +0.58%  "

Set-TextCell $ws "D43" "2.43"
$ws.Range("E43").Value = "  -10.29%  "

Set-TextCell $ws "B44" "InjectiveProtocol"
Set-TextCell $ws "C44" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell $ws "D44" "26.25"
$ws.Range("E44").Value = "  +4.66%  "

Set-TextCell $ws "B45" "USDe"
Set-TextCell $ws "C45" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws "D45" "1.00"
$ws.Range("E45").Value = "  -0.01%  "

Set-TextCell $ws "D46" "0.244"
$ws.Range("E46").Value = "  -3.69%  "

Set-TextCell $ws "B47" "Stellar"
Set-TextCell $ws "C47" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D47" "0.108"
$ws.Range("E47").Value = "  -0.99%  "

Set-TextCell $ws "B48" "Fetch.AI"
Set-TextCell $ws "C48" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell $ws "D48" "1.98"
$ws.Range("E48").Value = "  -3.11%  "

Set-TextCell $ws "D49" "117.95"
$ws.Range("E49").Value = "  -4.79%  "

$ws.Range("E50").Value = "  +6.99%  "

Set-TextCell $ws "D51" "0.0₃0490"
$ws.Range("E51").Value = "  -8.03%  "
